# Schedule generation now creates slightly better schedules.
#
# The "Red" train-schedule sheet previously listed 7 drivers/departure
# times (one per Employee, starting at Employee 1 / train 101) using a
# naive "6.7.0", "6.14.0", "6.21.0", ... increment-by-7 pattern.
#
# The regenerated schedule drops the Employee 1 / train 101 row entirely
# and recomputes the departure times for the remaining trains using a
# better pattern (two alternating shifts instead of one long runway).

$wb = $excel.ActiveWorkbook
$red = $wb.Worksheets.Item("Red")

# Drop the first train's row (Train 101 / Employee 1) from the schedule.
$red.Range("A2:C2").ClearContents()

# Recompute the departure times for the remaining trains (rows shift up
# logically, but stay on the same physical rows 3-8; the Driver names in
# column B are unchanged).
$red.Range("C3").Value = "06.00.00"
$red.Range("C4").Value = "14.00.00"
$red.Range("C5").Value = "6.7.0"
$red.Range("C6").Value = "2.7.0"
$red.Range("C7").Value = "6.14.0"
$red.Range("C8").Value = "2.14.0"
